$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(7)
$tr = $sh.TextFrame.TextRange

# Step 1: "strange man" -> "freaky man" (splits off a new, separately-styled run)
$full = $tr.Text
$idx = $full.IndexOf("strange man")
$len = "strange man".Length
$part = $tr.Characters($idx + 1, $len)
$part.Text = "freaky man"

# Step 2: drop "and physics " from "...teaching mathematics and physics in high school..."
# by rewriting the remainder of the text (from the closing paren onward) as a single run,
# so the split stays at exactly the two boundaries the edit introduced.
$full2 = $tr.Text
$idx2 = $full2.IndexOf(")")
$len2 = $full2.Length - $idx2
$tail = $tr.Characters($idx2 + 1, $len2)
$tail.Text = "). He has dedicated his life to studying mathematical sciences and loves everything about it. He has been teaching mathematics in high school for ten years."
